$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The review from rontiddler560@gmail.com (row 3) is being removed entirely;
# all rows below it shift up by one.
$ws.Rows.Item(3).Delete()

# Row deletion in this runtime does not re-anchor the Hyperlinks collection,
# so rebuild it from scratch to match the new (post-shift) layout.
$ws.Cells.Hyperlinks.Delete()

$links = @(
    @{ Cell = "C2";  Address = "mailto:leviadlevi22@gmail.com";     Display = "leviadlevi22@gmail.com" },
    @{ Cell = "D2";  Address = "mailto:gazittalia1@gmail.com";      Display = "gazittalia1@gmail.com" },
    @{ Cell = "C3";  Address = "mailto:gregneri12@gmail.com";       Display = "gregneri12@gmail.com" },
    @{ Cell = "C4";  Address = "mailto:snizzvered@gmail.com";       Display = "snizzvered@gmail.com" },
    @{ Cell = "C5";  Address = "mailto:budoyoni2@gmail.com";        Display = "budoyoni2@gmail.com" },
    @{ Cell = "C7";  Address = "mailto:hermanliran@gmail.com";      Display = "hermanliran@gmail.com" },
    @{ Cell = "C8";  Address = "mailto:gazittalia1@gmail.com";      Display = "gazittalia1@gmail.com" },
    @{ Cell = "D8";  Address = "mailto:hermanliran@gmail.com";      Display = "hermanliran@gmail.com" },
    @{ Cell = "C9";  Address = "mailto:freelancernachus@gmail.com"; Display = "freelancernachus@gmail.com" },
    @{ Cell = "C10"; Address = "mailto:nevilgreen@gmail.com";       Display = "nevilgreen@gmail.com" },
    @{ Cell = "D10"; Address = "mailto:vikicrestina@gmail.com";     Display = "vikicrestina@gmail.com" },
    @{ Cell = "C11"; Address = "mailto:veredsnir12@gmail.com";      Display = "veredsnir12@gmail.com" },
    @{ Cell = "D11"; Address = "mailto:kevinkors122@gmail.com";     Display = "kevinkors122@gmail.com" },
    @{ Cell = "C12"; Address = "mailto:stevewonder3001@gmail.com";  Display = "stevewonder3001@gmail.com" },
    @{ Cell = "D12"; Address = "mailto:budoyoni@gmail.com";         Display = "budoyoni@gmail.com" },
    @{ Cell = "C13"; Address = "mailto:stclerari834@gmail.com";     Display = "stclerari834@gmail.com" },
    @{ Cell = "C14"; Address = "mailto:stcydouel274@gmail.com";     Display = "stcydouel274@gmail.com" },
    @{ Cell = "C15"; Address = "mailto:kevinkors122@gmail.com";     Display = "kevinkors122@gmail.com" },
    @{ Cell = "D15"; Address = "mailto:sinuspai@gmail.com";         Display = "sinuspai@gmail.com" }
)

foreach ($link in $links) {
    $ws.Hyperlinks.Add($ws.Range($link.Cell), $link.Address, "", "", $link.Display)
}

$ws.Range("B3").Select()
